# Import ADJ column from excel budget file.
# Adds three new header columns (ADJ1, ADJ2, ADJ3) to the budget upload
# template right after the existing "MAR" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "ADJ1"
$ws.Range("T1").Value = "ADJ2"
$ws.Range("U1").Value = "ADJ3"

$ws.Range("S4").Select()
